$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Scrum_(Date)" to "Scrum_4_19_2022"
$ws.Name = "Scrum_4_19_2022"

# Fill in Jonathan Gamble's scrum update row (row 3)
$ws.Range("B3").Value = "Finished the notification and property suggestion UI elements"
$ws.Range("C3").Value = "Deliver the final deliverables for the project"
$ws.Range("D3").Value = "N/A"

# Row 3 grows taller to fit the new wrapped text
$ws.Rows.Item(3).RowHeight = 30

# Move the active selection to D4
$ws.Range("D4").Select()
